# Update header labels: D1 -> Jan_2026 (new), E1 -> Dec_2025, F1 -> Nov_2025 (Oct_2025 column dropped)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Jan_2026"
$ws.Range("E1").Value = "Dec_2025"
$ws.Range("F1").Value = "Nov_2025"

# Full replacement data set (quant engine refresh), sorted by Jan_2026 holding % descending
$data = @(
    ,@('INE090A01021', 'ICICI Bank Limited', 9.564712, 0, 0, 9.564712, 9.564712)
    ,@('INE040A01034', 'HDFC Bank Limited', 8.956498, 0, 0, 8.956498, 8.956498)
    ,@('INE423A01024', 'Adani Enterprises Limited', 8.873132999999999, 9.170431000000001, 4.216446, -0.2972980000000014, 4.656686999999999)
    ,@('INE202B01038', 'Piramal Finance Ltd', 8.520021, 9.898523000000001, 9.044435999999999, -1.378502000000001, -0.5244149999999994)
    ,@('INE406A01037', 'Aurobindo Pharma Limited', 6.533598, 5.966757, 5.997543, 0.5668409999999993, 0.5360549999999993)
    ,@('INE364U01010', 'Adani Green Energy Limited', 5.842646, 6.487632, 8.031215, -0.6449859999999994, -2.188568999999999)
    ,@('INE795G01014', 'HDFC Life Insurance Co Ltd', 5.796873, 5.543845, 5.477494, 0.2530279999999996, 0.3193789999999996)
    ,@('INE917I01010', 'Bajaj Auto Limited', 3.971147, 3.604163, 3.392926, 0.3669840000000004, 0.5782210000000001)
    ,@('INE237A01036', 'Kotak Mahindra Bank Limited', 3.67925, 0, 0, 3.67925, 3.67925)
    ,@('INE200M01039', 'Varun Beverages Limited', 2.929109, 0, 0, 2.929109, 2.929109)
    ,@('INE522F01014', 'Coal India Ltd', 2.211141, 0, 0, 2.211141, 2.211141)
    ,@('INE417T01026', 'PB Fintech Limited', 1.753087, 0, 0, 1.753087, 1.753087)
    ,@('INE075A01022', 'Wipro Ltd', 1.452308, 1.504776, 0, -0.05246799999999996, 1.452308)
    ,@('INE127D01025', 'HDFC Asset Management Company Ltd', 1.256642, 0, 0, 1.256642, 1.256642)
    ,@('INE259A01022', 'Colgate-Palmolive (India) Ltd', 0.897848, 0, 0, 0.897848, 0.897848)
    ,@('INE775A01035', 'Samvardhana Motherson International Ltd', 0.360054, 7.212751, 3.878306, -6.852697, -3.518252)
    ,@('INE280A01028', 'Titan Company Limited', 0, 2.93265, 0, -2.93265, 0)
    ,@('INE271C01023', 'DLF Limited', 0, 6.482575, 6.614819, -6.482575, -6.614819)
    ,@('INE245A01021', 'Tata Power Company Limited', 0, 3.680464, 5.204439, -3.680464, -5.204439)
    ,@('INE237A01028', 'Kotak Mahindra Bank Limited', 0, 3.701086, 3.462638, -3.701086, -3.462638)
    ,@('INE192A01025', 'Tata Consumer Products Ltd', 0, 1.687128, 1.608529, -1.687128, -1.608529)
    ,@('INE423A20016', 'Adani Enterprises Limited Rights', 0, 0, 0.098166, 0, -0.098166)
    ,@('INE467B01029', 'Tata Consultancy Services Limited', 0, 0.749743, 0, -0.749743, 0)
    ,@('INE0CZ201020', 'ANTHEM BIOSCIENCES LIMITED', 0, 0, 5.865315, 0, -5.865315)
    ,@('INE647A01010', 'SRF Limited', 0, 0, 2.683271, 0, -2.683271)
    ,@('INE094A01015', 'Hindustan Petroleum Corporation Ltd', 0, 2.976955, 0, -2.976955, 0)
    ,@('INE030A01027', 'Hindustan Unilever Limited', 0, 1.46324, 0, -1.46324, 0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $r = $data[$i]
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = "quant Momentum Fund"
    $ws.Cells.Item($row, 4).Value = $r[2]
    $ws.Cells.Item($row, 5).Value = $r[3]
    $ws.Cells.Item($row, 6).Value = $r[4]
    $ws.Cells.Item($row, 7).Value = $r[5]
    $ws.Cells.Item($row, 8).Value = $r[6]
}
